$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 311, shifting existing rows (311-359) down to (312-360)
$ws.Rows.Item(311).Insert()

# Populate the newly inserted row 311 with the new record's data
$ws.Cells.Item(311, 1).Value = 7
$ws.Cells.Item(311, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(311, 3).Value = "Ñuble"
$ws.Cells.Item(311, 4).Value = 45218
$ws.Cells.Item(311, 5).Value = 16
$ws.Cells.Item(311, 6).Value = 100112024
$ws.Cells.Item(311, 7).Value = "Choclo"
$ws.Cells.Item(311, 8).Value = "Dulce o Americano"
$ws.Cells.Item(311, 9).Value = "Primera"
$ws.Cells.Item(311, 10).Value = 50
$ws.Cells.Item(311, 11).Value = 40000
$ws.Cells.Item(311, 12).Value = 40000
$ws.Cells.Item(311, 13).Value = 40000
$ws.Cells.Item(311, 14).Value = "$/malla 70 unidades"
$ws.Cells.Item(311, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(311, 16).Value = 571
$ws.Cells.Item(311, 17).Value = 70
$ws.Cells.Item(311, 18).Value = "Hortaliza"
